$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(3, 1).Value = "SC-5,SC-5 (2),CM-6 b"
$ws.Cells.Item(4, 1).Value = "AC-6 (9),AC-6 (8),AU-12 (3),AU-7 a,CM-5 (1),AU-7 b,AU-8 b"
$ws.Cells.Item(5, 1).Value = "AC-17 (1),AC-17 (9),CM-7 b,CM-6 b"
$ws.Cells.Item(8, 1).Value = "IA-2 (12),IA-2 (11)"
$ws.Cells.Item(10, 1).Value = "CM-7 (5) (b),CM-7 (2)"
$ws.Cells.Item(12, 1).Value = "AC-7 b,AC-7 a"
$ws.Cells.Item(15, 1).Value = "IA-2,AU-3 (1),IA-8"
$ws.Cells.Item(16, 1).Value = "AC-6 (10),CM-6 b"
$ws.Cells.Item(17, 1).Value = "AU-12 a,AU-3,MA-4 (1) (a),AU-12 c,AU-3 (1)"
$ws.Cells.Item(19, 1).Value = "IA-5 (1) (b),IA-5 (1) (a),CM-6 b"
$ws.Cells.Item(21, 1).Value = "AC-12,MA-4 (7),SC-10,MA-4 e"
$ws.Cells.Item(22, 1).Value = "AU-12 a,AU-14 (1),AU-3,AU-7 (1),AU-7 a,MA-4 (1) (a),AU-6 (4),CM-6 b,CM-5 (1),AU-3 (1)"
$ws.Cells.Item(25, 1).Value = "AU-12 a,AU-3,MA-4 (1) (a),AU-12 c,AU-3 (1)"
$ws.Cells.Item(29, 1).Value = "SC-8 (1),SC-8 (2),SC-8"
$ws.Cells.Item(31, 1).Value = "AU-12 a,AU-3,MA-4 (1) (a),AC-2 (4),AU-12 c,AU-3 (1)"
$ws.Cells.Item(38, 1).Value = "SI-11 b,AU-9"
$ws.Cells.Item(39, 1).Value = "AU-3,CM-6 b"
$ws.Cells.Item(45, 1).Value = "AC-8 c 1, AC-8 c 2, AC-8 c 3,AC-8 b,AC-8 a"
$ws.Cells.Item(53, 1).Value = "MA-4 (6),SC-13"
$ws.Cells.Item(55, 1).Value = "AC-17 (2),SC-8"
$ws.Cells.Item(63, 1).Value = "AU-5 a,AU-5 (1)"
$ws.Cells.Item(65, 1).Value = "IA-2 (2),CM-6 b"
$ws.Cells.Item(67, 1).Value = "AU-12 a,AU-3,MA-4 (1) (a),AU-12 c,AU-3 (1)"
$ws.Cells.Item(69, 1).Value = "AU-12 a,AU-12 (3),AU-12 c,AU-7 a,CM-6 b,CM-5 (1),AU-7 b,AU-8 b"
$ws.Cells.Item(71, 1).Value = "AU-3,AU-4 (1)"
$ws.Cells.Item(77, 1).Value = "AU-12 a,AU-3,MA-4 (1) (a),AC-2 (4),AU-12 c,AU-3 (1)"
$ws.Cells.Item(80, 1).Value = "IA-2 (2),IA-2 (3),IA-2 (1),IA-2 (4)"
$ws.Cells.Item(81, 1).Value = "CM-6 b,CM-5 (3)"
$ws.Cells.Item(86, 1).Value = "AU-12 a,AU-3,MA-4 (1) (a),AU-12 c,AU-3 (1)"
$ws.Cells.Item(88, 1).Value = "CM-5 (1),AC-2 (4),AC-6 (9),AU-12 c"
$ws.Cells.Item(89, 1).Value = "IA-2 (2),IA-2,IA-2 (3),IA-2 (5),IA-2 (4)"
$ws.Cells.Item(90, 1).Value = "IA-2 (12),IA-2 (11)"
$ws.Cells.Item(96, 1).Value = "SC-8 (1),SC-8,AC-18 (1)"
$ws.Cells.Item(102, 1).Value = "AU-12 a,AU-3,MA-4 (1) (a),AU-12 c,AU-3 (1)"
$ws.Cells.Item(111, 1).Value = "AU-5 a,AU-5 b"
$ws.Cells.Item(119, 1).Value = "AU-12 a,AU-3,MA-4 (1) (a),AU-12 c,AU-3 (1)"
$ws.Cells.Item(124, 1).Value = "AU-12 a,AU-3,MA-4 (1) (a),AU-12 c,AU-3 (1)"
$ws.Cells.Item(128, 1).Value = "CM-7 a,IA-5 (1) (c),CM-6 b"
$ws.Cells.Item(136, 1).Value = "AC-11 (1),AC-11 b"
$ws.Cells.Item(139, 1).Value = "SI-6 d,SI-6 b,CM-3 (5)"
$ws.Cells.Item(148, 1).Value = "AU-12 a,AU-14 (1),AU-3,MA-4 (1) (a),AU-12 c,AU-3 (1)"
$ws.Cells.Item(157, 1).Value = "AU-12 a,AU-3,MA-4 (1) (a),AU-12 c,AU-3 (1)"
$ws.Cells.Item(159, 1).Value = "AC-17 (2),SC-8"
$ws.Cells.Item(181, 1).Value = "SC-3,CM-6 b"
